# Updates the cryptos list figures (price/volume columns) to reflect the
# latest scrape, and fixes the row order of Kaspa / RenderToken (rows 36-37)
# which had swapped places in the source ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '94.301.42'
$ws.Range("E2").Value = '  +2.35%  '
$ws.Range("D3").Value = '3.110.98'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.21'
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '616.35'
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  +4.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.391'
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.825'
$ws.Range("E10").Value = '  +12.85%  '
$ws.Range("D11").Value = '3.114.96'
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.199'
$ws.Range("E12").Value = '  -1.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000245'
$ws.Range("E13").Value = '  -2.80%  '
$ws.Range("D14").Value = '93.940.98'
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.65'
$ws.Range("E15").Value = '  +1.22%  '
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("D17").Value = '3.696.29'
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").Value = '3.159.11'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.72'
$ws.Range("E19").Value = '  +1.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.92'
$ws.Range("E20").Value = '  +0.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.94'
$ws.Range("E21").Value = '  +1.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '449.40'
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("E23").Value = '  -1.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.00'
$ws.Range("E24").Value = '  -3.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.29'
$ws.Range("E25").Value = '  +5.24%  '
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '86.47'
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("E28").Value = '  +3.35%  '
$ws.Range("D29").Value = '3.291.14'
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.260'
$ws.Range("E31").Value = '  +14.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.180'
$ws.Range("E32").Value = '  +8.02%  '
$ws.Range("E33").Value = '  -7.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.31'
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.94'
$ws.Range("E36").Value = '  -1.34%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.162'
$ws.Range("E37").Value = '  -3.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.17'
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.92'
$ws.Range("E39").Value = '  -1.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.454'
$ws.Range("E40").Value = '  +4.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '479.29'
$ws.Range("E41").Value = '  -0.70%  '
$ws.Range("E42").Value = '  +7.81%  '
$ws.Range("E43").Value = '  -1.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.73'
$ws.Range("E44").Value = '  -9.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.30'
$ws.Range("E45").Value = '  -4.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '160.40'
$ws.Range("E47").Value = '  +0.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.694'
$ws.Range("E48").Value = '  -0.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.87'
$ws.Range("E49").Value = '  -2.89%  '
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("E51").Value = '  -3.95%  '
